# Update countries & provincias Spain
# Applies the COVID data refresh captured in the diff:
#  - Update "Datos actualizados" timestamp text
#  - Update several countries' statistics (Casos totales, Nuevos casos,
#    Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes)
#  - Because the data is kept sorted by "Casos totales" descending, two
#    pairs of adjacent rows swap their country name / figures:
#      * Moldavia / Costa Rica (rows 65-66)
#      * Islas Malvinas / Montserrat (rows 214-215)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Header / timestamp text -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 24 de Agosto de 2020 a las 22:27"

# --- Simple numeric updates (country stays on the same row) -----------------

# Row 4: Estados Unidos
$ws.Range("B4").Value = 5906040
$ws.Range("C4").Value = 31894
$ws.Range("D4").Value = 3200228
$ws.Range("E4").Value = 2524873
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 335
$ws.Range("H4").Value = 180939

# Row 23: Alemania
$ws.Range("B23").Value = 236113
$ws.Range("C23").Value = 1624
$ws.Range("D23").Value = 209600
$ws.Range("E23").Value = 17177
$ws.Range("F23").Value = 0
$ws.Range("G23").Value = 4
$ws.Range("H23").Value = 9336

# Row 27: Canada
$ws.Range("B27").Value = 125120
$ws.Range("C27").Value = 224
$ws.Range("D27").Value = 111210
$ws.Range("E27").Value = 4832
$ws.Range("F27").Value = 0
$ws.Range("G27").Value = 5
$ws.Range("H27").Value = 9078

# Row 33: Israel
$ws.Range("B33").Value = 104472
$ws.Range("C33").Value = 1809
$ws.Range("D33").Value = 81642
$ws.Range("E33").Value = 21983
$ws.Range("F33").Value = 0
$ws.Range("G33").Value = 13
$ws.Range("H33").Value = 847

# Row 127: Tunez
$ws.Range("B127").Value = 2893
$ws.Range("C127").Value = 75
$ws.Range("D127").Value = 1454
$ws.Range("E127").Value = 1368

# Row 158: Principado de Andorra
$ws.Range("B158").Value = 1060
$ws.Range("C158").Value = 15
$ws.Range("D158").Value = 877
$ws.Range("E158").Value = 130

# --- Rows whose country swapped position because of the new sort order ------

# Rows 65-66: Costa Rica overtakes Moldavia
$ws.Range("A65").Value = "Costa Rica"
$ws.Range("B65").Value = 34463
$ws.Range("C65").Value = 643
$ws.Range("D65").Value = 12758
$ws.Range("E65").Value = 21343
$ws.Range("F65").Value = 0
$ws.Range("G65").Value = 7
$ws.Range("H65").Value = 362

$ws.Range("A66").Value = "Moldavia"
$ws.Range("B66").Value = 33828
$ws.Range("C66").Value = 350
$ws.Range("D66").Value = 23570
$ws.Range("E66").Value = 9313
$ws.Range("F66").Value = 0
$ws.Range("G66").Value = 5
$ws.Range("H66").Value = 945

# Rows 214-215: Montserrat overtakes Islas Malvinas
$ws.Range("A214").Value = "Montserrat"
$ws.Range("B214").Value = 13
$ws.Range("C214").Value = 0
$ws.Range("D214").Value = 12
$ws.Range("E214").Value = 0
$ws.Range("F214").Value = 0
$ws.Range("G214").Value = 0
$ws.Range("H214").Value = 1

$ws.Range("A215").Value = "Islas Malvinas"
$ws.Range("B215").Value = 13
$ws.Range("C215").Value = 0
$ws.Range("D215").Value = 13
$ws.Range("E215").Value = 0
$ws.Range("F215").Value = 0
$ws.Range("G215").Value = 0
$ws.Range("H215").Value = 0
